# Auto commit update: refresh Metrics figures and selections.
$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update raw metric values (B2:B13) ---
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 402758.63000000006
$metrics.Range("B3").Value  = 344658.97
$metrics.Range("B4").Value  = 123261.01
$metrics.Range("B5").Value  = 16272
$metrics.Range("B6").Value  = 5605465.7400000002
$metrics.Range("B7").Value  = 4745011.9300000006
$metrics.Range("B8").Value  = 1655217.8900000004
$metrics.Range("B9").Value  = 218979
$metrics.Range("B10").Value = 34070846.730000004
$metrics.Range("B11").Value = 32020287.090000004
$metrics.Range("B12").Value = 11936939.929999994
$metrics.Range("B13").Value = 1316609

# Move the Metrics sheet's active-cell selection from D14 to D18.
$metrics.Activate()
$metrics.Range("D18").Select()

# --- today sheet: move the active-cell selection from D6 to E8 ---
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("E8").Select()
